$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The runs/balls/fours/sixes columns (C,D,E,F) store numeric-looking values
# as text, so prefix with an apostrophe to keep them text-typed instead of
# letting Excel auto-convert them to numbers.

# Row 2
$ws.Range("C2").Value = "'0"
$ws.Range("D2").Value = "'0"
$ws.Range("E2").Value = "'0"
$ws.Range("F2").Value = "'0"

# Row 3
$ws.Range("C3").Value = "'15"
$ws.Range("D3").Value = "'15"
$ws.Range("E3").Value = "'2"
$ws.Range("F3").Value = "'0"

# Row 4
$ws.Range("C4").Value = "'12"
$ws.Range("D4").Value = "'7"
$ws.Range("E4").Value = "'0"
$ws.Range("F4").Value = "'1"

# Row 5
$ws.Range("C5").Value = "'15"
$ws.Range("D5").Value = "'7"
$ws.Range("E5").Value = "'1"
$ws.Range("F5").Value = "'1"

# Row 7
$ws.Range("C7").Value = "'3"
$ws.Range("D7").Value = "'7"
$ws.Range("E7").Value = "'0"
$ws.Range("F7").Value = "'0"

# Row 8
$ws.Range("C8").Value = "'9"
$ws.Range("D8").Value = "'10"
$ws.Range("E8").Value = "'1"
$ws.Range("F8").Value = "'0"
